# "Generate Report for Handoff" — refresh the latest-handoff timestamps for
# the abee8654-13f7-4858-9d05-2aebe7a81458.md row (row 7 on each sheet) to
# reflect a newly generated handoff report.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-16 02:38:13"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-16 02:38:18"

# Overview sheet: column G = "Latest HO Xliff Generate Date" (tracks the
# most recent handoff across locales, i.e. the de-de value above)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-16 02:38:18"
